$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 410.75
$ws.Range("I12").Value = 410.75
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 410.75
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -240.75
$ws.Range("N12").ClearContents() | Out-Null
$ws.Range("H17").Value = 936.86487
$ws.Range("J17").Value = 942.971
$ws.Range("L17").Value = 2828.913
$ws.Range("N17").Value = -3164.913
$ws.Range("H19").Value = 4536.6665
$ws.Range("I19").Value = 6116.6665
$ws.Range("J19").Value = 2166.6667
$ws.Range("K19").Value = 6116.6665
$ws.Range("L19").Value = 2166.6667
$ws.Range("M19").Value = -5941.6665
$ws.Range("N19").Value = -2516.6667
$ws.Range("H116").Value = 3626.2703
$ws.Range("I116").Value = 3951.9546
$ws.Range("K116").Value = 3951.9546
$ws.Range("M116").Value = -509.9546
$ws.Range("H132").Value = 2649.3914
$ws.Range("I132").Value = 1798.0646
$ws.Range("K132").Value = 5394.1938
$ws.Range("M132").Value = -2864.1938
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1537.5
$ws.Range("I2").Value = 1260
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1260
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1147
$ws.Range("N2").Value = -2226
$ws.Range("H88").Value = 2514.875
$ws.Range("I88").Value = 2732.1428
$ws.Range("J88").Value = 2210.7
$ws.Range("K88").Value = 2732.1428
$ws.Range("L88").Value = 2210.7
$ws.Range("M88").Value = -2326.1428
$ws.Range("N88").Value = -3022.7
$ws.Range("H91").Value = 2514.875
$ws.Range("I91").Value = 2732.1428
$ws.Range("J91").Value = 2210.7
$ws.Range("K91").Value = 2732.1428
$ws.Range("L91").Value = 2210.7
$ws.Range("M91").Value = -1328.1428
$ws.Range("N91").Value = -5018.7
$ws.Range("H116").Value = 1537.5
$ws.Range("I116").Value = 1260
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1260
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1034
$ws.Range("N116").Value = -6588
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1537.5
$ws.Range("I3").Value = 1260
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1260
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1146
$ws.Range("N3").Value = -2228
$ws.Range("H86").Value = 7436.6
$ws.Range("I86").Value = 8958
$ws.Range("J86").Value = 3886.6667
$ws.Range("K86").Value = 8958
$ws.Range("L86").Value = 3886.6667
$ws.Range("M86").Value = -7835
$ws.Range("N86").Value = -6132.6667
$ws.Range("H89").Value = 7436.6
$ws.Range("I89").Value = 8958
$ws.Range("J89").Value = 3886.6667
$ws.Range("K89").Value = 44790
$ws.Range("L89").Value = 19433.3335
$ws.Range("M89").Value = -39174
$ws.Range("N89").Value = -30665.3335
$ws.Range("H134").Value = 2447.6616
$ws.Range("I134").Value = 2033.96
$ws.Range("J134").Value = 3826.6667
$ws.Range("K134").Value = 6101.88
$ws.Range("L134").Value = 11480.0001
$ws.Range("M134").Value = -3566.88
$ws.Range("N134").Value = -16550.0001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 776
$ws.Range("I16").Value = 784.2857
$ws.Range("J16").Value = 756.6667
$ws.Range("K16").Value = 784.2857
$ws.Range("L16").Value = 756.6667
$ws.Range("M16").Value = -497.2857
$ws.Range("N16").Value = -1330.6667
$ws.Range("H31").Value = 2757.6843
$ws.Range("I31").Value = 1683.5135
$ws.Range("J31").Value = 4744.9
$ws.Range("K31").Value = 1683.5135
$ws.Range("L31").Value = 4744.9
$ws.Range("M31").Value = -1388.5135
$ws.Range("N31").Value = -5334.9
$ws.Range("H34").Value = 2757.6843
$ws.Range("I34").Value = 1683.5135
$ws.Range("J34").Value = 4744.9
$ws.Range("K34").Value = 1683.5135
$ws.Range("L34").Value = 4744.9
$ws.Range("M34").Value = -1481.5135
$ws.Range("N34").Value = -5148.9
$ws.Range("H58").Value = 3137.82
$ws.Range("I58").Value = 3370.4473
$ws.Range("K58").Value = 3370.4473
$ws.Range("M58").Value = -3167.4473
$ws.Range("H94").Value = 9227.700000000001
$ws.Range("I94").Value = 1600
$ws.Range("J94").Value = 10075.223
$ws.Range("K94").Value = 1600
$ws.Range("L94").Value = 10075.223
$ws.Range("M94").Value = -1149
$ws.Range("N94").Value = -10977.223
$ws.Range("H113").Value = 776
$ws.Range("I113").Value = 784.2857
$ws.Range("J113").Value = 756.6667
$ws.Range("K113").Value = 784.2857
$ws.Range("L113").Value = 756.6667
$ws.Range("M113").Value = 1385.7143
$ws.Range("N113").Value = -5096.6667
$ws.Range("H122").Value = 1534.3889
$ws.Range("I122").Value = 877.5
$ws.Range("J122").Value = 2059.9
$ws.Range("K122").Value = 2632.5
$ws.Range("L122").Value = 6179.700000000001
$ws.Range("M122").Value = -182.5
$ws.Range("N122").Value = -11079.7
$ws.Range("H132").Value = 2674.1785
$ws.Range("I132").Value = 1225.3334
$ws.Range("J132").Value = 4345.923
$ws.Range("K132").Value = 3676.0002
$ws.Range("L132").Value = 13037.769
$ws.Range("M132").Value = -1146.0002
$ws.Range("N132").Value = -18097.769
$ws.Range("H134").Value = 1523.2632
$ws.Range("I134").Value = 976.7895
$ws.Range("J134").Value = 2616.2104
$ws.Range("K134").Value = 2930.3685
$ws.Range("L134").Value = 7848.6312
$ws.Range("M134").Value = -395.3685
$ws.Range("N134").Value = -12918.6312
$ws.Range("H136").Value = 3137.82
$ws.Range("I136").Value = 3370.4473
$ws.Range("K136").Value = 10111.3419
$ws.Range("M136").Value = -7561.341899999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6750
$ws.Range("I3").Value = 8500
$ws.Range("K3").Value = 25500
$ws.Range("M3").Value = -25388
$ws.Range("H129").Value = 2630.2693
$ws.Range("J129").Value = 2113
$ws.Range("L129").Value = 6339
$ws.Range("N129").Value = -16339
$ws.Range("H131").Value = 2164.9844
$ws.Range("J131").Value = 1778.6604
$ws.Range("L131").Value = 5335.9812
$ws.Range("N131").Value = -15415.9812
$ws.Range("H133").Value = 11795.125
$ws.Range("I133").Value = 7515
$ws.Range("J133").Value = 13221.833
$ws.Range("K133").Value = 22545
$ws.Range("L133").Value = 39665.499
$ws.Range("M133").Value = -17485
$ws.Range("N133").Value = -49785.499
$ws.Range("H134").Value = 5802.294
$ws.Range("I134").Value = 5962.636
$ws.Range("J134").Value = 5508.3335
$ws.Range("K134").Value = 17887.908
$ws.Range("L134").Value = 16525.0005
$ws.Range("M134").Value = -12817.908
$ws.Range("N134").Value = -26665.0005
$ws.Range("H136").Value = 3612.2222
$ws.Range("I136").Value = 3438.75
$ws.Range("K136").Value = 10316.25
$ws.Range("M136").Value = -5216.25
$ws.Range("H139").Value = 923.3333
$ws.Range("I139").Value = 923.3333
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 2769.9999
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 2370.0001
$ws.Range("N139").ClearContents() | Out-Null
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3930.9424
$ws.Range("I70").Value = 3978.9473
$ws.Range("J70").Value = 3800.6428
$ws.Range("K70").Value = 3978.9473
$ws.Range("L70").Value = 3800.6428
$ws.Range("M70").Value = -3708.9473
$ws.Range("N70").Value = -4340.6428
$ws.Range("H73").Value = 3930.9424
$ws.Range("I73").Value = 3978.9473
$ws.Range("J73").Value = 3800.6428
$ws.Range("K73").Value = 3978.9473
$ws.Range("L73").Value = 3800.6428
$ws.Range("M73").Value = -3042.9473
$ws.Range("N73").Value = -5672.6428
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8299.333000000001
$ws.Range("I132").Value = 2815.84
$ws.Range("J132").Value = 18091.285
$ws.Range("K132").Value = 8447.52
$ws.Range("L132").Value = 54273.855
$ws.Range("M132").Value = -5917.52
$ws.Range("N132").Value = -59333.855
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1671.463
$ws.Range("I132").Value = 1056.4865
$ws.Range("K132").Value = 3169.4595
$ws.Range("M132").Value = -639.4594999999999
$ws.Range("H136").Value = 18055590
$ws.Range("I136").Value = 23280632
$ws.Range("J136").Value = 772753.4399999999
$ws.Range("K136").Value = 69841896
$ws.Range("L136").Value = 2318260.32
$ws.Range("M136").Value = -69839346
